$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("micsss")

# Delete the entire "ASSAY CATEGORY" row (row 4) - shifts everything below up by one.
$ws.Rows.Item(4).Delete()

# Select the row that now occupies row 4 (matches Excel's post-delete selection behavior).
$ws.Range("A4:XFD4").Select()
